# New .ttl from Google sheet has been generated
# Applies the refreshed export to the "vocabulary" sheet:
#  - ConceptScheme / PREFIX URIs point at the new purl path
#  - placeholder/test metadata values are cleared out
#  - a stray duplicate "dct:creator" metadata row is removed (rows shift up)
#  - the now-blank "vars:test" / "vars:computerscientist" example rows are
#    reset back to empty template rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated ConceptScheme URI / PREFIX URI
$ws.Range("B1").Value = "http://purl.org/m4m-dk-Test4/variables/"
$ws.Range("C3").Value = "http://purl.org/m4m-dk-Test4/variables/"

# Clear out placeholder test values from the metadata block
$ws.Range("B10").Value = ""
$ws.Range("B11").Value = ""
$ws.Range("B12").Value = ""

# Remove the duplicate "dct:creator" row (A13:T13, value "Minka") -
# remaining metadata/header/data rows shift up by one
$ws.Rows(13).Delete()

# The example term rows (now at 19 and 20 after the shift) were only ever
# placeholders ("vars:test"/"test" and "vars:computerscientist"/...) -
# reset them back to blank "vars:" template rows
$ws.Range("A19").Value = "vars:"
$ws.Range("B19").Value = ""

$ws.Range("A20").Value = "vars:"
$ws.Range("B20").Value = ""
$ws.Range("E20").Value = ""
